# Update countries & provincias Spain
# Applies the refreshed COVID numbers and the resulting re-ranking of a
# handful of countries whose "Casos totales" overtook their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($r, $name, $vals) {
    if ($name -ne $null) {
        $ws.Cells.Item($r, 1).Value = $name
    }
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
    $ws.Cells.Item($r, 8).Value = $vals[6]
}

# --- Header: "Datos actualizados a ..." timestamp -----------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 00:10"

# --- Worldwide total row (Estados Unidos, row 4) -------------------------
Set-Row 4 $null @(1766220, 20417, 497960, 1164994, 0, 1159, 103266)

# --- Guyana overtakes Brunei and Islas Caimanes (rows 164-166) ----------
Set-Row 164 "Guyana"         @(150, 11, 67, 72, 0, 0, 11)
Set-Row 165 "Brunei"         @(141, 0, 138, 1, 0, 0, 2)
Set-Row 166 "Islas Caimanes" @(140, 0, 67, 72, 0, 0, 1)

# --- Curazao overtakes Fiyi (rows 198-199) -------------------------------
Set-Row 198 "Curazao" @(18, 0, 14, 3, 0, 0, 1)
Set-Row 199 "Fiyi"     @(18, 0, 15, 3, 0, 0, 0)

# --- Santa Lucia overtakes Belice (rows 200-201) -------------------------
Set-Row 200 "Santa Lucia" @(18, 0, 18, 0, 0, 0, 0)
Set-Row 201 "Belice"      @(18, 0, 16, 0, 0, 0, 2)

# --- Seychelles overtakes Montserrat (rows 210-211) -----------------------
Set-Row 210 "Seychelles" @(11, 0, 11, 0, 0, 0, 0)
Set-Row 211 "Montserrat" @(11, 0, 10, 0, 0, 0, 1)

# --- Bonaire, San Eustaquio y Saba overtakes San Bartolome (rows 215-216) -
Set-Row 215 "Bonaire, San Eustaquio y Saba" @(6, 0, 6, 0, 0, 0, 0)
Set-Row 216 "San Bartolome"                 @(6, 0, 6, 0, 0, 0, 0)
